$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = [double]"29.223446"
$ws.Range("H2").Value = [double]"87.670338"
$ws.Range("I2").Value = [double]"0.0169041244192178"
$ws.Range("J2").Value = [double]"0.0169041244192178"
$ws.Range("M2").Value = [double]"1.819857"
$ws.Range("N2").Value = [double]"5.459571"
$ws.Range("O2").Value = [double]"0.01485317462584607"
$ws.Range("P2").Value = [double]"0.01485317462584607"
$ws.Range("Q2").Value = [double]"53.182492767222"
$ws.Range("R2").Value = [double]"478.642434904998"
$ws.Range("S2").Value = [double]"0.0002510799118956707"
$ws.Range("T2").Value = [double]"0.0002510799118956708"
$ws.Range("G3").Value = [double]"29.223446"
$ws.Range("H3").Value = [double]"87.670338"
$ws.Range("I3").Value = [double]"0.0169041244192178"
$ws.Range("J3").Value = [double]"0.0169041244192178"
$ws.Range("O3").Value = [double]"0.726618572334523"
$ws.Range("P3").Value = [double]"0.7266185723345231"
$ws.Range("Q3").Value = [double]"2601.69209217176"
$ws.Range("R3").Value = [double]"23415.22882954584"
$ws.Range("S3").Value = [double]"0.01228285075205719"
$ws.Range("T3").Value = [double]"0.01228285075205719"
$ws.Range("G4").Value = [double]"29.223446"
$ws.Range("H4").Value = [double]"87.670338"
$ws.Range("I4").Value = [double]"0.0169041244192178"
$ws.Range("J4").Value = [double]"0.0169041244192178"
$ws.Range("M4").Value = [double]"31.52924033333333"
$ws.Range("N4").Value = [double]"94.58772099999999"
$ws.Range("O4").Value = [double]"0.257333028084772"
$ws.Range("P4").Value = [double]"0.257333028084772"
$ws.Range("Q4").Value = [double]"921.3930523021886"
$ws.Range("R4").Value = [double]"8292.537470719697"
$ws.Range("S4").Value = [double]"0.004349989523919056"
$ws.Range("T4").Value = [double]"0.004349989523919056"
$ws.Range("G5").Value = [double]"29.223446"
$ws.Range("H5").Value = [double]"87.670338"
$ws.Range("I5").Value = [double]"0.0169041244192178"
$ws.Range("J5").Value = [double]"0.0169041244192178"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.1464426666666667"
$ws.Range("N5").Value = [double]"0.439328"
$ws.Range("O5").Value = [double]"0.001195224954858853"
$ws.Range("P5").Value = [double]"0.001195224954858853"
$ws.Range("Q5").Value = [double]"4.279559361429333"
$ws.Range("R5").Value = [double]"38.516034252864"
$ws.Range("S5").Value = [double]"2.020423134588803E-05"
$ws.Range("T5").Value = [double]"2.020423134588803E-05"
$ws.Range("I6").Value = [double]"0.9471112884046843"
$ws.Range("J6").Value = [double]"0.9471112884046842"
$ws.Range("M6").Value = [double]"1.819857"
$ws.Range("N6").Value = [double]"5.459571"
$ws.Range("O6").Value = [double]"0.01485317462584607"
$ws.Range("P6").Value = [double]"0.01485317462584607"
$ws.Range("Q6").Value = [double]"2979.73074476857"
$ws.Range("R6").Value = [double]"26817.57670291713"
$ws.Range("S6").Value = [double]"0.01406760935678483"
$ws.Range("T6").Value = [double]"0.01406760935678483"
$ws.Range("I7").Value = [double]"0.9471112884046843"
$ws.Range("J7").Value = [double]"0.9471112884046842"
$ws.Range("O7").Value = [double]"0.726618572334523"
$ws.Range("P7").Value = [double]"0.7266185723345231"
$ws.Range("S7").Value = [double]"0.6881886522225223"
$ws.Range("T7").Value = [double]"0.6881886522225223"
$ws.Range("I8").Value = [double]"0.9471112884046843"
$ws.Range("J8").Value = [double]"0.9471112884046842"
$ws.Range("M8").Value = [double]"31.52924033333333"
$ws.Range("N8").Value = [double]"94.58772099999999"
$ws.Range("O8").Value = [double]"0.257333028084772"
$ws.Range("P8").Value = [double]"0.257333028084772"
$ws.Range("Q8").Value = [double]"51624.19178014017"
$ws.Range("R8").Value = [double]"464617.7260212615"
$ws.Range("S8").Value = [double]"0.2437230157784473"
$ws.Range("T8").Value = [double]"0.2437230157784472"
$ws.Range("I9").Value = [double]"0.9471112884046843"
$ws.Range("J9").Value = [double]"0.9471112884046842"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.1464426666666667"
$ws.Range("N9").Value = [double]"0.439328"
$ws.Range("O9").Value = [double]"0.001195224954858853"
$ws.Range("P9").Value = [double]"0.001195224954858853"
$ws.Range("Q9").Value = [double]"239.7769254466489"
$ws.Range("R9").Value = [double]"2157.99232901984"
$ws.Range("S9").Value = [double]"0.001132011046929798"
$ws.Range("T9").Value = [double]"0.001132011046929798"
$ws.Range("G10").Value = [double]"37.39212666666667"
$ws.Range("H10").Value = [double]"112.17638"
$ws.Range("I10").Value = [double]"0.02162924801792661"
$ws.Range("J10").Value = [double]"0.0216292480179266"
$ws.Range("M10").Value = [double]"1.819857"
$ws.Range("N10").Value = [double]"5.459571"
$ws.Range("O10").Value = [double]"0.01485317462584607"
$ws.Range("P10").Value = [double]"0.01485317462584607"
$ws.Range("Q10").Value = [double]"68.04832345922"
$ws.Range("R10").Value = [double]"612.4349111329801"
$ws.Range("S10").Value = [double]"0.0003212629978359988"
$ws.Range("T10").Value = [double]"0.0003212629978359988"
$ws.Range("G11").Value = [double]"37.39212666666667"
$ws.Range("H11").Value = [double]"112.17638"
$ws.Range("I11").Value = [double]"0.02162924801792661"
$ws.Range("J11").Value = [double]"0.0216292480179266"
$ws.Range("O11").Value = [double]"0.726618572334523"
$ws.Range("P11").Value = [double]"0.7266185723345231"
$ws.Range("Q11").Value = [double]"3328.929800344267"
$ws.Range("R11").Value = [double]"29960.3682030984"
$ws.Range("S11").Value = [double]"0.01571621331545514"
$ws.Range("T11").Value = [double]"0.01571621331545514"
$ws.Range("G12").Value = [double]"37.39212666666667"
$ws.Range("H12").Value = [double]"112.17638"
$ws.Range("I12").Value = [double]"0.02162924801792661"
$ws.Range("J12").Value = [double]"0.0216292480179266"
$ws.Range("M12").Value = [double]"31.52924033333333"
$ws.Range("N12").Value = [double]"94.58772099999999"
$ws.Range("O12").Value = [double]"0.257333028084772"
$ws.Range("P12").Value = [double]"0.257333028084772"
$ws.Range("Q12").Value = [double]"1178.945348247775"
$ws.Range("R12").Value = [double]"10610.50813422998"
$ws.Range("S12").Value = [double]"0.005565919887649608"
$ws.Range("T12").Value = [double]"0.005565919887649607"
$ws.Range("G13").Value = [double]"37.39212666666667"
$ws.Range("H13").Value = [double]"112.17638"
$ws.Range("I13").Value = [double]"0.02162924801792661"
$ws.Range("J13").Value = [double]"0.0216292480179266"
$ws.Range("K13").Value = [double]"1"
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.1464426666666667"
$ws.Range("N13").Value = [double]"0.439328"
$ws.Range("O13").Value = [double]"0.001195224954858853"
$ws.Range("P13").Value = [double]"0.001195224954858853"
$ws.Range("Q13").Value = [double]"5.475802741404445"
$ws.Range("R13").Value = [double]"49.28222467264001"
$ws.Range("S13").Value = [double]"2.585181698585726E-05"
$ws.Range("T13").Value = [double]"2.585181698585726E-05"
$ws.Range("G14").Value = [double]"24.817167"
$ws.Range("H14").Value = [double]"74.45150100000001"
$ws.Range("I14").Value = [double]"0.01435533915817136"
$ws.Range("J14").Value = [double]"0.01435533915817136"
$ws.Range("M14").Value = [double]"1.819857"
$ws.Range("N14").Value = [double]"5.459571"
$ws.Range("O14").Value = [double]"0.01485317462584607"
$ws.Range("P14").Value = [double]"0.01485317462584607"
$ws.Range("Q14").Value = [double]"45.163695085119"
$ws.Range("R14").Value = [double]"406.4732557660711"
$ws.Range("S14").Value = [double]"0.0002132223593295653"
$ws.Range("T14").Value = [double]"0.0002132223593295653"
$ws.Range("G15").Value = [double]"24.817167"
$ws.Range("H15").Value = [double]"74.45150100000001"
$ws.Range("I15").Value = [double]"0.01435533915817136"
$ws.Range("J15").Value = [double]"0.01435533915817136"
$ws.Range("O15").Value = [double]"0.726618572334523"
$ws.Range("P15").Value = [double]"0.7266185723345231"
$ws.Range("Q15").Value = [double]"2209.41182412252"
$ws.Range("R15").Value = [double]"19884.70641710268"
$ws.Range("S15").Value = [double]"0.01043085604448835"
$ws.Range("T15").Value = [double]"0.01043085604448835"
$ws.Range("G16").Value = [double]"24.817167"
$ws.Range("H16").Value = [double]"74.45150100000001"
$ws.Range("I16").Value = [double]"0.01435533915817136"
$ws.Range("J16").Value = [double]"0.01435533915817136"
$ws.Range("M16").Value = [double]"31.52924033333333"
$ws.Range("N16").Value = [double]"94.58772099999999"
$ws.Range("O16").Value = [double]"0.257333028084772"
$ws.Range("P16").Value = [double]"0.257333028084772"
$ws.Range("Q16").Value = [double]"782.466422735469"
$ws.Range("R16").Value = [double]"7042.197804619221"
$ws.Range("S16").Value = [double]"0.003694102894756139"
$ws.Range("T16").Value = [double]"0.003694102894756138"
$ws.Range("G17").Value = [double]"24.817167"
$ws.Range("H17").Value = [double]"74.45150100000001"
$ws.Range("I17").Value = [double]"0.01435533915817136"
$ws.Range("J17").Value = [double]"0.01435533915817136"
$ws.Range("K17").Value = [double]"1"
$ws.Range("L17").Value = [double]"0.3333333333333333"
$ws.Range("M17").Value = [double]"0.1464426666666667"
$ws.Range("N17").Value = [double]"0.439328"
$ws.Range("O17").Value = [double]"0.001195224954858853"
$ws.Range("P17").Value = [double]"0.001195224954858853"
$ws.Range("Q17").Value = [double]"3.634292114592"
$ws.Range("R17").Value = [double]"32.708629031328"
$ws.Range("S17").Value = [double]"1.715785959730888E-05"
$ws.Range("T17").Value = [double]"1.715785959730888E-05"
